# Insert a new "Skill Description" column (full skill names) right after
# the existing SkillCode column (A), shifting the former SFIA Level /
# Keycode / Description columns one place to the right (B->C, C->D, D->E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new column.
$ws.Columns("B").Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Skill Description"

# Map SkillCode (column A) -> full skill name for the new column B.
$skillNames = @{
    "Autonomy"  = "Autonomy";
    "Influence" = "Influence";
    "Complexity" = "Complexity";
    "Knowledge" = "Knowledge";
    "SWDN" = "Software design";
    "PROG" = "Programming/software development";
    "TEST" = "Testing";
    "DTAN" = "Data modelling and design";
    "DENG" = "Data engineering";
    "MADE" = "MADE";
    "REQM" = "Requirements definition and management";
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 32) { $lastRow = 32 }

for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value2
    if ($code -ne $null -and $code -ne "") {
        $name = $skillNames[$code]
        if ($name -eq $null) { $name = $code }
        $ws.Cells.Item($r, 2).Value = $name
    }
}
